$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the "matrix size 7500" table (column D, rows 2-21) with the new
#    Opencl measurements.
# ---------------------------------------------------------------------------
$colD1 = @(
    0.49505399999999999,
    0.49049399999999999,
    0.48364099999999999,
    0.39807799999999999,
    0.40021600000000002,
    0.383185,
    0.40973799999999999,
    0.37506200000000001,
    0.38073499999999999,
    0.37915199999999999,
    0.37596600000000002,
    0.38241999999999998,
    0.384432,
    0.383183,
    0.37625900000000001,
    0.37137399999999998,
    0.37882199999999999,
    0.39555400000000002,
    0.37514900000000001,
    0.38949699999999998
)
for ($i = 0; $i -lt $colD1.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $colD1[$i]
}

# ---------------------------------------------------------------------------
# 2. Update the "12 Threads" table (column D, rows 24-38) with the new
#    Opencl measurements.
# ---------------------------------------------------------------------------
$colD2 = @(
    0.124069,
    0.12784899999999999,
    0.13491800000000001,
    0.14202500000000001,
    0.148592,
    0.179539,
    0.19530500000000001,
    0.22969400000000001,
    0.25361400000000001,
    0.28656900000000002,
    0.31492399999999998,
    0.35635699999999998,
    0.39061000000000001,
    0.43513800000000002,
    0.481097
)
for ($i = 0; $i -lt $colD2.Count; $i++) {
    $row = 24 + $i
    $ws.Cells.Item($row, 4).Value = $colD2[$i]
}

# ---------------------------------------------------------------------------
# 3. Convert both charts from 3-D clustered column charts to plain 2-D
#    clustered column charts, and adjust the first chart's value-axis
#    scaling.
# ---------------------------------------------------------------------------
$chart1 = $ws.ChartObjects(1).Chart
$chart1.ChartType = 51
$valAx1 = $chart1.Axes(2)
$valAx1.MaximumScale = 1.4
$valAx1.MinimumScale = 0

$chart2 = $ws.ChartObjects(2).Chart
$chart2.ChartType = 51

# ---------------------------------------------------------------------------
# 4. Move the sheet selection to D24:D38 (active cell D24), matching where
#    the user was working when the edit was made.
# ---------------------------------------------------------------------------
$ws.Range("D24:D38").Select()
